$d = $word.ActiveDocument

# Update the estimated cost figure in the project scope table:
# R$ 12.483,20 -> R$ 14.979,84
$d.Content.Find.ClearFormatting()
$d.Content.Find.Execute(
    "12.483,20",  # FindText
    $true,        # MatchCase
    $false,       # MatchWholeWord
    $false,       # MatchWildcards
    $false,       # MatchSoundsLike
    $false,       # MatchAllWordForms
    $true,        # Forward
    1,            # Wrap (wdFindContinue)
    $false,       # Format
    "14.979,84",  # ReplaceWith
    2             # Replace (wdReplaceAll)
)
